{"js": "// Fix the spelling of the surname \"Nieuwenhuys\" -> \"Nieuwenhuijs\"\n// (the author's real, officially registered surname is \"Nieuwenhuijs\",\n// i.e. spelled with \"ij\" instead of \"y\").\nconst body = context.document.body;\n\nconst results = body.search(\"Nieuwenhuys\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (const range of results.items) {\n  range.insertText(\"Nieuwenhuijs\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix the spelling of the surname \"Nieuwenhuys\" -> \"Nieuwenhuijs\"\n# (the author's real, officially registered surname is \"Nieuwenhuijs\",\n# i.e. spelled with \"ij\" instead of \"y\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Nieuwenhuys\"\n$find.Replacement.Text = \"Nieuwenhuijs\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdReplaceAll = 2\n$find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $null, \"Nieuwenhuijs\", 2)\n"}
